# "07 - 6000 to 7000.xlsx" — re-sort the data table by Profit (Last 5 Year
# Avg.), column B, descending, and highlight the row that lands on row 9
# (• Finolex Ind) with a yellow fill, then leave that row selected —
# mirroring the manual "select row -> sort -> highlight" edit captured in
# the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-sort the data rows (A2:F14) by column B, descending -----------
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("B1"), 0, 2)
$sortObj.SetRange($ws.Range("A2:F14"))
$sortObj.Header = 2
$sortObj.Apply()

# --- Highlight the row that is now row 9 with a solid yellow fill -----
$ws.Range("A9:F9").Interior.Color = 65535

# --- Leave the whole row selected, like the saved sheetView shows -----
$ws.Range("A9:XFD9").Select()
